# Updated symbol list on Sat Jan 14 08:51:16 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (col D) and "Volume(1h)" (col E) columns for the
# crypto rows whose quotes moved since the last scrape. All of these cells
# are stored as plain text (e.g. "313.47", "9.17%") rather than numbers, so
# we briefly force Text number format before writing the values (otherwise
# Excel would auto-convert them to numeric/percentage values) and restore
# the original "Normal" style afterwards so no visible formatting changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Contiguous blocks covering every cell touched below. Kept as separate
# (non-union) ranges because only the first area of a disjoint/union range
# reliably picks up the NumberFormat change.
$r1 = $ws.Range("D2:E28")
$r1.NumberFormat = "@"
$r2 = $ws.Range("D40:E45")
$r2.NumberFormat = "@"
$r3 = $ws.Range("E46:E47")
$r3.NumberFormat = "@"

$ws.Range("D2").Value = "313.47"
$ws.Range("E2").Value = "9.17%"
$ws.Range("D3").Value = "32.44"
$ws.Range("E3").Value = "9.66%"
$ws.Range("D4").Value = "5.333"
$ws.Range("E4").Value = "4.25%"
$ws.Range("D5").Value = "0.07686"
$ws.Range("E5").Value = "14.78%"
$ws.Range("D6").Value = "7.880"
$ws.Range("E6").Value = "7.27%"
$ws.Range("D7").Value = "3.728"
$ws.Range("E7").Value = "9.53%"
$ws.Range("D8").Value = "1.628"
$ws.Range("E8").Value = "19.39%"
$ws.Range("D9").Value = "0.9193"
$ws.Range("E9").Value = "0.74%"
$ws.Range("D10").Value = "0.01769"
$ws.Range("E10").Value = "2,635.37%"
$ws.Range("D11").Value = "0.1727"
$ws.Range("E11").Value = "8.87%"
$ws.Range("D12").Value = "0.07590"
$ws.Range("E12").Value = "12.49%"
$ws.Range("D13").Value = "0.08247"
$ws.Range("E13").Value = "6.95%"
$ws.Range("D14").Value = "0.03034"
$ws.Range("E14").Value = "3.62%"
$ws.Range("D15").Value = "0.09891"
$ws.Range("E15").Value = "10.16%"
$ws.Range("D16").Value = "0.001521"
$ws.Range("E16").Value = "-3.12%"
$ws.Range("D17").Value = "0.04565"
$ws.Range("E17").Value = "1.56%"
$ws.Range("D18").Value = "0.006096"
$ws.Range("E18").Value = "-2.45%"
$ws.Range("D19").Value = "3.474"
$ws.Range("E19").Value = "0.95%"
$ws.Range("D20").Value = "2.246"
$ws.Range("E20").Value = "1.01%"
$ws.Range("D21").Value = "0.3324"
$ws.Range("E21").Value = "3.45%"
$ws.Range("D22").Value = "0.1335"
$ws.Range("E22").Value = "1.99%"
$ws.Range("D23").Value = "4.234"
$ws.Range("E23").Value = "4.26%"
$ws.Range("D24").Value = "0.1624"
$ws.Range("E24").Value = "2.68%"
$ws.Range("D25").Value = "0.001220"
$ws.Range("E25").Value = "2.26%"
$ws.Range("D26").Value = "0.004511"
$ws.Range("E26").Value = "9.51%"
$ws.Range("D27").Value = "0.0001297"
$ws.Range("E27").Value = "8.17%"
$ws.Range("E28").Value = "7.45%"
$ws.Range("D40").Value = "0.04651"
$ws.Range("E40").Value = "9.25%"
$ws.Range("D41").Value = "0.007196"
$ws.Range("E41").Value = "6.11%"
$ws.Range("D42").Value = "0.1374"
$ws.Range("E42").Value = "10.84%"
$ws.Range("D43").Value = "0.002255"
$ws.Range("E43").Value = "4.01%"
$ws.Range("D44").Value = "0.01459"
$ws.Range("E44").Value = "9.82%"
$ws.Range("D45").Value = "0.00006203"
$ws.Range("E45").Value = "8.85%"
$ws.Range("E46").Value = "-3.83%"
$ws.Range("E47").Value = "-0.68%"

# Restore original (no explicit number format) appearance.
$r1.Style = "Normal"
$r2.Style = "Normal"
$r3.Style = "Normal"
